$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells (B:E data, matching the original inlineStr cell types)
# so numeric-looking strings like "1.005" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.242.01"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.914.85"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("D5").Value = "0.7368"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "243.94"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").Value = "0.3120"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("D9").Value = "26.97"
$ws.Range("E9").Value = "  -4.15%  "
$ws.Range("D10").Value = "0.06936"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "0.07975"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "0.7695"
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("D13").Value = "1.939.55"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "5.272"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "91.24"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "30.229.38"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "5.833"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "244.40"
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("D20").Value = "0.000007822"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "2.163.51"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "6.646"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "9.354"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "164.88"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "18.89"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "0.1266"
$ws.Range("E28").Value = "  -5.99%  "
$ws.Range("D29").Value = "2.121"
$ws.Range("E29").Value = "  -9.12%  "
$ws.Range("D30").Value = "1.354"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "1.548"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "4.335"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "4.058"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").Value = "0.05160"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "1.286"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "0.7447"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "2.780"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.01932"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "2.793"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "6.335"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "75.53"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").Value = "0.4443"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "1.934"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "1.004"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "0.8361"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "101.19"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.600"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.837"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.077.58"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "37.08"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1202"
$ws.Range("E51").Value = "  +2.63%  "
